$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.676.94'
$ws.Range('E2').Value = '  +1.66%  '
$ws.Range('D3').Value = '1.635.55'
$ws.Range('E3').Value = '  +1.92%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.94'
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('E6').Value = '  +1.99%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  +1.13%  '
$ws.Range('E9').Value = '  +1.72%  '
$ws.Range('E10').Value = '  +2.91%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0838'
$ws.Range('E11').Value = '  +2.88%  '
$ws.Range('D12').Value = '1.863.42'
$ws.Range('E12').Value = '  +1.94%  '
$ws.Range('D13').Value = '1.641.09'
$ws.Range('E13').Value = '  +2.27%  '
$ws.Range('E14').Value = '  +1.28%  '
$ws.Range('E15').Value = '  +2.43%  '
$ws.Range('D16').Value = '26.683.73'
$ws.Range('E16').Value = '  +1.84%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.01'
$ws.Range('E17').Value = '  +1.65%  '
$ws.Range('D18').Value = '0.0₃0741'
$ws.Range('E18').Value = '  +1.81%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '208.81'
$ws.Range('E19').Value = '  +3.84%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.00'
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.30'
$ws.Range('E21').Value = '  +0.70%  '
$ws.Range('E22').Value = '  +0.89%  '
$ws.Range('E23').Value = '  +2.73%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.90'
$ws.Range('E24').Value = '  +1.54%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.36'
$ws.Range('E25').Value = '  +1.69%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  -0.80%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.73'
$ws.Range('E28').Value = '  +2.57%  '
$ws.Range('E29').Value = '  +1.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0520'
$ws.Range('E30').Value = '  +5.92%  '
$ws.Range('E31').Value = '  -0.10%  '
$ws.Range('E32').Value = '  +1.09%  '
$ws.Range('E33').Value = '  +0.70%  '
$ws.Range('E34').Value = '  +1.89%  '
$ws.Range('E35').Value = '  +0.66%  '
$ws.Range('D36').Value = '1.168.59'
$ws.Range('E36').Value = '  +0.54%  '
$ws.Range('E37').Value = '  -1.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.808'
$ws.Range('E38').Value = '  +2.85%  '
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.503'
$ws.Range('E40').Value = '  +1.43%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.32'
$ws.Range('E41').Value = '  +0.32%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.793'
$ws.Range('E42').Value = '  +1.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.38'
$ws.Range('E43').Value = '  +0.92%  '
$ws.Range('D44').Value = '1.775.21'
$ws.Range('E44').Value = '  +2.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.32'
$ws.Range('E45').Value = '  +0.68%  '
$ws.Range('E46').Value = '  +1.39%  '
$ws.Range('E47').Value = '  -1.67%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '54.65'
$ws.Range('E48').Value = '  +0.97%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.409'
$ws.Range('E50').Value = '  +0.57%  '
$ws.Range('E51').Value = '  +4.33%  '
